{"js": "// Update the two-digit multiplication answers in the table.\n// Each old equation string is unique in the document, so a direct\n// search-and-replace of the full \"A\u00d7B=C\" text for each cell is safe.\nconst pairs = [\n  [\"76\u00d786=6536\", \"24\u00d756=1344\"],\n  [\"37\u00d777=2849\", \"94\u00d780=7520\"],\n  [\"59\u00d732=1888\", \"21\u00d749=1029\"],\n  [\"50\u00d767=3350\", \"52\u00d719=988\"],\n  [\"39\u00d712=468\", \"76\u00d719=1444\"],\n  [\"74\u00d774=5476\", \"69\u00d756=3864\"],\n  [\"30\u00d723=690\", \"59\u00d739=2301\"],\n  [\"64\u00d760=3840\", \"70\u00d739=2730\"],\n  [\"98\u00d784=8232\", \"93\u00d744=4092\"],\n  [\"96\u00d794=9024\", \"70\u00d721=1470\"],\n  [\"64\u00d785=5440\", \"13\u00d747=611\"],\n  [\"57\u00d762=3534\", \"79\u00d769=5451\"],\n  [\"84\u00d726=2184\", \"98\u00d760=5880\"],\n  [\"69\u00d783=5727\", \"53\u00d729=1537\"],\n  [\"12\u00d795=1140\", \"27\u00d797=2619\"],\n  [\"77\u00d771=5467\", \"39\u00d747=1833\"],\n  [\"30\u00d732=960\", \"76\u00d714=1064\"],\n  [\"31\u00d771=2201\", \"37\u00d721=777\"],\n  [\"84\u00d766=5544\", \"42\u00d729=1218\"],\n  [\"42\u00d764=2688\", \"66\u00d776=5016\"],\n  [\"25\u00d799=2475\", \"94\u00d764=6016\"],\n  [\"54\u00d789=4806\", \"43\u00d720=860\"],\n  [\"17\u00d719=323\", \"30\u00d725=750\"],\n  [\"44\u00d742=1848\", \"46\u00d721=966\"],\n  [\"86\u00d736=3096\", \"32\u00d743=1376\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication answers in the table.\n# Each old equation string is unique in the document, so Find/Replace\n# of the full \"A\u00d7B=C\" text for each cell is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"76\u00d786=6536\", \"24\u00d756=1344\"),\n    @(\"37\u00d777=2849\", \"94\u00d780=7520\"),\n    @(\"59\u00d732=1888\", \"21\u00d749=1029\"),\n    @(\"50\u00d767=3350\", \"52\u00d719=988\"),\n    @(\"39\u00d712=468\", \"76\u00d719=1444\"),\n    @(\"74\u00d774=5476\", \"69\u00d756=3864\"),\n    @(\"30\u00d723=690\", \"59\u00d739=2301\"),\n    @(\"64\u00d760=3840\", \"70\u00d739=2730\"),\n    @(\"98\u00d784=8232\", \"93\u00d744=4092\"),\n    @(\"96\u00d794=9024\", \"70\u00d721=1470\"),\n    @(\"64\u00d785=5440\", \"13\u00d747=611\"),\n    @(\"57\u00d762=3534\", \"79\u00d769=5451\"),\n    @(\"84\u00d726=2184\", \"98\u00d760=5880\"),\n    @(\"69\u00d783=5727\", \"53\u00d729=1537\"),\n    @(\"12\u00d795=1140\", \"27\u00d797=2619\"),\n    @(\"77\u00d771=5467\", \"39\u00d747=1833\"),\n    @(\"30\u00d732=960\", \"76\u00d714=1064\"),\n    @(\"31\u00d771=2201\", \"37\u00d721=777\"),\n    @(\"84\u00d766=5544\", \"42\u00d729=1218\"),\n    @(\"42\u00d764=2688\", \"66\u00d776=5016\"),\n    @(\"25\u00d799=2475\", \"94\u00d764=6016\"),\n    @(\"54\u00d789=4806\", \"43\u00d720=860\"),\n    @(\"17\u00d719=323\", \"30\u00d725=750\"),\n    @(\"44\u00d742=1848\", \"46\u00d721=966\"),\n    @(\"86\u00d736=3096\", \"32\u00d743=1376\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
